# Move the "Performance results" slide earlier in the deck so that it
# appears right after the two intro "Writing the tiled kernel" slides,
# motivating the step-by-step kernel-writing work that follows.
#
# Before: ... Step1, Step2, Step3, Step3a, Step3b, Step4, Performance results
# After:  ... Performance results, Step1, Step2, Step3, Step3a, Step3b, Step4

$p = $ppt.ActivePresentation

# Locate the "Performance results" slide by its title text rather than a
# hard-coded index, so the script is resilient to slide numbering.
$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $text = $shape.TextFrame.TextRange.Text.Trim()
            if ($text -eq "Performance results") {
                $targetIndex = $i
            }
        }
    }
    if ($targetIndex -ne -1) {
        break
    }
}

# Move it to slot 9 - right after the two earlier "Writing the tiled
# kernel" slides, and before what used to be the "Step 1" opening slide.
$destIndex = 9

if ($targetIndex -ne -1 -and $targetIndex -ne $destIndex) {
    $s = $p.Slides.Item($targetIndex)
    $s.MoveTo($destIndex)
}
